$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply formatting to the whole new block (rows 1082-1111) to match existing data rows
$newRange = $ws.Range("A1082:G1111")
$newRange.Font.Name = "Times New Roman"
$newRange.Font.Size = 13
$newRange.EntireRow.RowHeight = 17

# Populate the 30 new rows for Unit 37
$ws.Range("A1082").Value = 'U37_01'
$ws.Range("B1082").Value = 37
$ws.Range("C1082").Value = 'Khu dân cư'
$ws.Range("D1082").Value = 'Neighborhood'
$ws.Range("E1082").Value = 'There are no burglars in the neighborhood'
$ws.Range("F1082").Value = 'In the neighborhood / trong khu phố'
$ws.Range("G1082").Value = 'N'

$ws.Range("A1083").Value = 'U37_02'
$ws.Range("B1083").Value = 37
$ws.Range("C1083").Value = 'Đi lang thang'
$ws.Range("D1083").Value = 'Wander'
$ws.Range("E1083").Value = 'I enjoy wandering around the city center'
$ws.Range("F1083").Value = 'wander around / đi lang thang xung quanh'
$ws.Range("G1083").Value = 'V'

$ws.Range("A1084").Value = 'U37_03'
$ws.Range("B1084").Value = 37
$ws.Range("C1084").Value = 'Sân chơi'
$ws.Range("D1084").Value = 'Playground'
$ws.Range("E1084").Value = 'Students can play at the school playground'
$ws.Range("F1084").Value = 'a school playground / sân trường'
$ws.Range("G1084").Value = 'N'

$ws.Range("A1085").Value = 'U37_04'
$ws.Range("B1085").Value = 37
$ws.Range("C1085").Value = 'Nhịp'
$ws.Range("D1085").Value = 'Pace'
$ws.Range("E1085").Value = 'She is learning at a steady pace'
$ws.Range("F1085").Value = 'at a steady pace / ở một nhịp độ ổn định'
$ws.Range("G1085").Value = 'N'

$ws.Range("A1086").Value = 'U37_05'
$ws.Range("B1086").Value = 37
$ws.Range("C1086").Value = 'Người đi đường'
$ws.Range("D1086").Value = 'Pedestrian'
$ws.Range("E1086").Value = 'We can walk on the pedestrian crossing'
$ws.Range("F1086").Value = 'a pedestrian crossing / vạch kẻ đường cho người đi bộ'
$ws.Range("G1086").Value = 'N'

$ws.Range("A1087").Value = 'U37_06'
$ws.Range("B1087").Value = 37
$ws.Range("C1087").Value = 'Vỉa hè'
$ws.Range("D1087").Value = 'Sidewalk'
$ws.Range("E1087").Value = 'At night, the city has busy sidewalks'
$ws.Range("F1087").Value = 'a busy sidewalk'
$ws.Range("G1087").Value = 'N'

$ws.Range("A1088").Value = 'U37_07'
$ws.Range("B1088").Value = 37
$ws.Range("C1088").Value = 'Làn đường'
$ws.Range("D1088").Value = 'Lane'
$ws.Range("E1088").Value = 'We drove in the right traffic lane.'
$ws.Range("F1088").Value = 'a traffic lane / làn đường giao thông'
$ws.Range("G1088").Value = 'N'

$ws.Range("A1089").Value = 'U37_08'
$ws.Range("B1089").Value = 37
$ws.Range("C1089").Value = 'Giao thông'
$ws.Range("D1089").Value = 'Traffic'
$ws.Range("E1089").Value = 'I want to avoid the rush-hour traffic'
$ws.Range("F1089").Value = 'rush-hour traffic / giao thông giờ cao điểm'
$ws.Range("G1089").Value = 'N'

$ws.Range("A1090").Value = 'U37_09'
$ws.Range("B1090").Value = 37
$ws.Range("C1090").Value = 'Rạp chiếu'
$ws.Range("D1090").Value = 'Theater'
$ws.Range("E1090").Value = 'There is a movie theater around the corner'
$ws.Range("F1090").Value = 'a movie theater'
$ws.Range("G1090").Value = 'N'

$ws.Range("A1091").Value = 'U37_10'
$ws.Range("B1091").Value = 37
$ws.Range("C1091").Value = 'Cửa hàng sách'
$ws.Range("D1091").Value = 'Bookstore'
$ws.Range("E1091").Value = 'Nowadays, we buy from online bookstores.'
$ws.Range("F1091").Value = 'An online bookstore / hiệu sách trực tuyến'
$ws.Range("G1091").Value = 'N'

$ws.Range("A1092").Value = 'U37_11'
$ws.Range("B1092").Value = 37
$ws.Range("C1092").Value = 'Tạp hóa'
$ws.Range("D1092").Value = 'Grocery'
$ws.Range("E1092").Value = 'Our grocery bill is too high'
$ws.Range("F1092").Value = 'The grocery bill'
$ws.Range("G1092").Value = 'N'

$ws.Range("A1093").Value = 'U37_12'
$ws.Range("B1093").Value = 37
$ws.Range("C1093").Value = 'Ngõ'
$ws.Range("D1093").Value = 'Alley'
$ws.Range("E1093").Value = 'My car cannot fit in the narrow alley'
$ws.Range("F1093").Value = 'A narrow alley'
$ws.Range("G1093").Value = 'N'

$ws.Range("A1094").Value = 'U37_13'
$ws.Range("B1094").Value = 37
$ws.Range("C1094").Value = 'Tòa nhà'
$ws.Range("D1094").Value = 'Building'
$ws.Range("E1094").Value = 'The mausoleum is a historic building'
$ws.Range("F1094").Value = 'a historic building / công trình lịch sử'
$ws.Range("G1094").Value = 'N'

$ws.Range("A1095").Value = 'U37_14'
$ws.Range("B1095").Value = 37
$ws.Range("C1095").Value = 'Xuống cấp'
$ws.Range("D1095").Value = 'Run-down'
$ws.Range("E1095").Value = 'We can fix this run-down building'
$ws.Range("F1095").Value = 'a run-down building'
$ws.Range("G1095").Value = 'Adj'

$ws.Range("A1096").Value = 'U37_15'
$ws.Range("B1096").Value = 37
$ws.Range("C1096").Value = 'Tượng trưng cho'
$ws.Range("D1096").Value = 'Represent'
$ws.Range("E1096").Value = 'Roses represent remance and love'
$ws.Range("F1096").Value = 'represent somebody or something / tượng trưng cho ai đó hoặc cái gì đó'
$ws.Range("G1096").Value = 'V'

$ws.Range("A1097").Value = 'U37_16'
$ws.Range("B1097").Value = 37
$ws.Range("C1097").Value = 'Tháp'
$ws.Range("D1097").Value = 'Tower'
$ws.Range("E1097").Value = 'The eiffel tower is in Paris, France'
$ws.Range("F1097").Value = 'The eiffel tower'
$ws.Range("G1097").Value = 'N'

$ws.Range("A1098").Value = 'U37_17'
$ws.Range("B1098").Value = 37
$ws.Range("C1098").Value = 'Giao thông vận tải'
$ws.Range("D1098").Value = 'Transportation'
$ws.Range("E1098").Value = 'What is your means of transportation to work every day'
$ws.Range("F1098").Value = 'means of transportation / phương tiện giao thông'
$ws.Range("G1098").Value = 'N'

$ws.Range("A1099").Value = 'U37_18'
$ws.Range("B1099").Value = 37
$ws.Range("C1099").Value = 'Dễ tiếp cận'
$ws.Range("D1099").Value = 'Accessible'
$ws.Range("E1099").Value = 'The library is accessible to everyone'
$ws.Range("F1099").Value = 'accessible to somebody / Ai đó có thể tiếp cận được'
$ws.Range("G1099").Value = 'Adj'

$ws.Range("A1100").Value = 'U37_19'
$ws.Range("B1100").Value = 37
$ws.Range("C1100").Value = 'Cảnh sát'
$ws.Range("D1100").Value = 'Police'
$ws.Range("E1100").Value = 'There is a police car outside'
$ws.Range("F1100").Value = 'a police car '
$ws.Range("G1100").Value = 'N'

$ws.Range("A1101").Value = 'U37_20'
$ws.Range("B1101").Value = 37
$ws.Range("C1101").Value = 'Dắt đi tham quan'
$ws.Range("D1101").Value = 'Show around'
$ws.Range("E1101").Value = 'I can show you around town sometime'
$ws.Range("F1101").Value = 'show somebody around town / dắt ai đi tham quan thị trấn'
$ws.Range("G1101").Value = 'V'

$ws.Range("A1102").Value = 'U37_21'
$ws.Range("B1102").Value = 37
$ws.Range("C1102").Value = 'Tài sản'
$ws.Range("D1102").Value = 'Property'
$ws.Range("E1102").Value = 'This painting is government property'
$ws.Range("F1102").Value = 'government property / tài sản chính phủ'
$ws.Range("G1102").Value = 'N'

$ws.Range("A1103").Value = 'U37_22'
$ws.Range("B1103").Value = 37
$ws.Range("C1103").Value = 'Đông đúc'
$ws.Range("D1103").Value = 'Crowded'
$ws.Range("E1103").Value = 'We avoid going to crowded areas on the weekend'
$ws.Range("F1103").Value = 'crowded areas'
$ws.Range("G1103").Value = 'Adj'

$ws.Range("A1104").Value = 'U37_23'
$ws.Range("B1104").Value = 37
$ws.Range("C1104").Value = 'Tiếng ồn'
$ws.Range("D1104").Value = 'Noise'
$ws.Range("E1104").Value = 'Children tend to make noises'
$ws.Range("F1104").Value = 'make a noise'
$ws.Range("G1104").Value = 'N'

$ws.Range("A1105").Value = 'U37_24'
$ws.Range("B1105").Value = 37
$ws.Range("C1105").Value = 'chuyển vào'
$ws.Range("D1105").Value = 'Move in'
$ws.Range("E1105").Value = 'The family moved in a new apartment '
$ws.Range("F1105").Value = 'move in an apartment '
$ws.Range("G1105").Value = 'V'

$ws.Range("A1106").Value = 'U37_25'
$ws.Range("B1106").Value = 37
$ws.Range("C1106").Value = 'Trung tâm'
$ws.Range("D1106").Value = 'Central'
$ws.Range("E1106").Value = 'Let''s meet at a central location'
$ws.Range("F1106").Value = 'a central location'
$ws.Range("G1106").Value = 'Adj'

$ws.Range("A1107").Value = 'U37_26'
$ws.Range("B1107").Value = 37
$ws.Range("C1107").Value = 'Chuỗi'
$ws.Range("D1107").Value = 'Chain'
$ws.Range("E1107").Value = 'A chain of events happened throughout tonight'
$ws.Range("F1107").Value = 'a chain of events'
$ws.Range("G1107").Value = 'N'

$ws.Range("A1108").Value = 'U37_27'
$ws.Range("B1108").Value = 37
$ws.Range("C1108").Value = 'Đa dạng'
$ws.Range("D1108").Value = 'Diverse'
$ws.Range("E1108").Value = 'Our country has a diverse culture'
$ws.Range("F1108").Value = 'A diverse culture'
$ws.Range("G1108").Value = 'Adj'

$ws.Range("A1109").Value = 'U37_28'
$ws.Range("B1109").Value = 37
$ws.Range("C1109").Value = 'Ngẫu nhiên'
$ws.Range("D1109").Value = 'Random'
$ws.Range("E1109").Value = 'We met on such a random occasion'
$ws.Range("F1109").Value = 'A random accasion / một dịp ngẫu nhiên'
$ws.Range("G1109").Value = 'Adj'

$ws.Range("A1110").Value = 'U37_29'
$ws.Range("B1110").Value = 37
$ws.Range("C1110").Value = 'Dày đặc'
$ws.Range("D1110").Value = 'Dense'
$ws.Range("E1110").Value = 'She can''t pass through the dense crowd'
$ws.Range("F1110").Value = 'a dense crowd'
$ws.Range("G1110").Value = 'Adj'

$ws.Range("A1111").Value = 'U37_30'
$ws.Range("B1111").Value = 37
$ws.Range("C1111").Value = 'Dưới mặt đất'
$ws.Range("D1111").Value = 'Underground'
$ws.Range("E1111").Value = 'There is an underground basement here'
$ws.Range("F1111").Value = 'an underground basement / một tầng hầm dưới mặt đất'
$ws.Range("G1111").Value = 'Adj'

# Widen column D to fit the new, longer English words (e.g. "Transportation")
$ws.Columns.Item(4).ColumnWidth = 14.8

# Move the selection/view to the end of the newly added data, as in the saved workbook
$ws.Range("G1111").Select()